$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C7").Value = -12.9086
$ws.Range("B9").Value = 5.424999999999998
$ws.Range("C12").Value = -11.1258
$ws.Range("B18").Value = 6.455899999999998
$ws.Range("B20").Value = 9.275099999999995
$ws.Range("C26").Value = -13.21560000000001
$ws.Range("B27").Value = 6.549300000000003
$ws.Range("C27").Value = -12.9157
$ws.Range("C29").Value = -11.29280000000001
$ws.Range("C37").Value = -13.47100000000001
$ws.Range("C38").Value = -12.9192
$ws.Range("C51").Value = -11.5561
$ws.Range("C55").Value = -13.54139999999999
$ws.Range("B69").Value = 5.829299999999998
$ws.Range("C69").Value = -12.5463
$ws.Range("C70").Value = -11.2967
$ws.Range("B76").Value = 5.547900000000001
$ws.Range("B82").Value = 5.644499999999998
$ws.Range("C83").Value = -14.2065
$ws.Range("C102").Value = -13.17000000000001
